# [Fonds de solidarite] Fix missing 2020-06-25 data
# Update nombre_aides (col C) and montant_total (col D) for the rows that
# were missing the 2020-06-25 data refresh. Values are stored as text
# (inline/shared strings), so each cell is forced to Text format before the
# write and the style is then reset to "Normal" to avoid leaving a stray
# number-format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "C2" "150"
Set-TextValue "D2" "330000.00"

Set-TextValue "C3" "802"
Set-TextValue "D3" "2010761.80"

Set-TextValue "C6" "16"
Set-TextValue "D6" "78500.00"

Set-TextValue "C9" "30"
Set-TextValue "D9" "66000.00"

Set-TextValue "C10" "243"
Set-TextValue "D10" "565990.66"

Set-TextValue "C11" "102"
Set-TextValue "D11" "314891.77"

Set-TextValue "C12" "23"
Set-TextValue "D12" "94000.00"

Set-TextValue "C21" "45"
Set-TextValue "D21" "111500.00"

Set-TextValue "C23" "100"
Set-TextValue "D23" "332600.00"

Set-TextValue "C25" "15"
Set-TextValue "D25" "31500.00"

Set-TextValue "C33" "419"
Set-TextValue "D33" "1026811.79"

Set-TextValue "C72" "734"
Set-TextValue "D72" "1892946.83"

Set-TextValue "C74" "86"
Set-TextValue "D74" "335000.00"

Set-TextValue "C83" "77"
Set-TextValue "D83" "173200.00"

Set-TextValue "C84" "338"
Set-TextValue "D84" "814972.09"

Set-TextValue "C85" "133"
Set-TextValue "D85" "414192.00"

Set-TextValue "C86" "38"
Set-TextValue "D86" "141709.01"

Set-TextValue "C87" "7"
Set-TextValue "D87" "27500.00"

$wb.Save()
